$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").ClearContents()
$ws.Range("H4").Select()
